# Update cryptocurrency price/volume figures on Sheet1.
# Values in column D ("Price") that would otherwise be re-interpreted by
# Excel as numbers (and so lose formatting such as trailing zeros or
# thousand-group separators) are entered with a leading apostrophe so
# that they remain plain text, matching how the source data is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "35.648.79"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "1.984.26"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'242.53"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.638"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").Value = "'56.99"
$ws.Range("E7").Value = "  +8.28%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'60.06"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").Value = "'0.925"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "'14.11"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "2.275.97"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "'5.26"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "1.989.15"
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "'17.26"
$ws.Range("E18").Value = "  +5.61%  "
$ws.Range("D19").Value = "35.530.35"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").Value = "'70.46"
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'233.29"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  +10.56%  "
$ws.Range("D27").Value = "'163.63"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "'19.50"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'4.80"
$ws.Range("E32").Value = "  -5.32%  "
$ws.Range("D33").Value = "'0.0588"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "'0.0897"
$ws.Range("E34").Value = "  +9.22%  "
$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  -6.46%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "'4.93"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").Value = "'0.0891"
$ws.Range("E44").Value = "  -4.36%  "
$ws.Range("D45").Value = "'91.00"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "1.379.06"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'7.43"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'15.44"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").Value = "'45.73"
$ws.Range("E51").Value = "  +1.84%  "

